$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The save_data regeneration now derives strikeouts from the "K" column
# (previously populated from a "Strike#" source). Column G holds these
# recalculated K values; write the new per-game values below.
$kValues = [ordered]@{
    2 = 1
    3 = 2
    4 = 1
    5 = 1
    6 = 1
    8 = 0
    9 = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 0
    34 = 1
    35 = 2
    36 = 0
    37 = 1
    38 = 1
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    43 = 2
    44 = 1
    45 = 0
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    50 = 2
    51 = 1
    52 = 1
    53 = 1
    54 = 2
    55 = 1
    56 = 0
    57 = 2
    58 = 1
    59 = 1
    60 = 0
    62 = 1
    63 = 0
    64 = 2
    65 = 1
    66 = 3
    67 = 1
    68 = 1
    69 = 1
    70 = 1
    71 = 2
    72 = 1
    73 = 2
    74 = 1
    75 = 2
    76 = 2
    77 = 1
    78 = 0
    79 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}

